# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-26 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 4
    3  = 7
    4  = 2
    5  = 5
    6  = 4
    7  = 5
    8  = 2
    9  = 6
    10 = 7
    11 = 3
    12 = 5
    13 = 4
    14 = 2
    15 = 2
    16 = 5
    17 = 8
    18 = 3
    19 = 4
    20 = 2
    21 = 3
    22 = 5
    23 = 6
    24 = 4
    25 = 2
    26 = 2
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
